$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.942.46"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "2.052.42"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.81%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "2.355.55"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "2.062.98"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "37.865.36"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.04%  "
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.47%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "1.500.70"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0919"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("E47").Value = "  +15.74%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "2.244.19"
$ws.Range("E51").Value = "  +1.47%  "
